$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Defects" column (I) was actually mislabeled - it held a last-update
# timestamp, not defect data. Rename the header to reflect reality.
$ws.Range("I1").Value = "Date of Last Update"

# Insert two fresh columns where the real "TicketID" duplicate and the
# real "Defects" data belong (pushing "Engineering Issues" / "Location of
# Issue" one column further right, from J/K to L/M).
$ws.Columns("J:K").Insert()

# New column J: duplicate of TicketID (column A) for each row.
$ws.Range("J1").Value = "TicketID"
$ws.Range("A2").Copy()
$ws.Range("J2").PasteSpecial(-4163)
$ws.Range("A3").Copy()
$ws.Range("J3").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# New column K: the real "Defects" data.
$ws.Range("K1").Value = "Defects"
$ws.Range("K2").Value = "Changes requested / Suggestion"
$ws.Range("K3").Value = "Damaged / Broken / Defective"
